$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Der er enighed om behovet for klare og konsekvente regler for tildeling af dansk statsborgerskab og bekymringer over ændringer af betingelserne efter ansøgninger, hvilket underminerer tilliden til systemet. Diskussionen om statsborgerskab omfatter også spørgsmål om diskrimination baseret på oprindelse, økonomiske kriterier og beskyttelse af grundlæggende frihedsrettigheder. Flere partier ønsker at sikre, at borgere, der opfylder kravene, behandles retfærdigt og ikke bliver påvirket af politiske skift eller fordomme.'
$ws.Range("B3").Value = 'At opnå dansk statsborgerskab kræver opfyldelse af strenge krav, herunder sprogkundskaber, selvforsørgelse og kendskab til danske værdier, hvilket anerkendes som et stort privilegium. Fejringen af de nye statsborgere understreger vigtigheden af deres bidrag til samfundet og muligheden for at deltage aktivt i demokratiet. Der er en generel enighed om, at statsborgerskab fremmer integration og tilhørsforhold, samtidig med at der er ønsker om at lette betingelserne for at gøre det lettere for flere at blive en del af det danske fællesskab.'

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'Der anerkendes udfordringer forbundet med islamisme og sharialov i visse ghettoer, samt problemer relateret til ikkevestlig indvandring.'
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Der er behov for, at politikere træffer beslutninger baseret på deres overbevisninger for at sikre Danmarks fremtid og velfærd. Regeringen anerkender forskellige synspunkter og har oprettet en ekspertgruppe for at forbedre borgerbehandlingen, men der er bekymringer om, at sindelagssamtaler kan overtræde retsstatsprincipperne. Der appelleres til ministeren om at sikre, at disse principper respekteres i statsborgerskabsansøgninger.'
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'Indlægget kritiserer Hr. Jesper Langballes udtalelser som perfide og faktuelt mangelfulde, samtidig med at der påpeges, at direkte tiltale ikke må anvendes i debatten.'
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 'Der er behov for at nedsætte en ekspertgruppe for at opnå resultater, men der er usikkerhed om tidshorisonten. Der søges klarhed over Venstres støtte til den tidligere beslutning om en undersøgelseskommission.'
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 'Der afholdes en række taler og bemærkninger fra forskellige ordførere i Folketinget, hvor der udveksles synspunkter om udlændingepolitik, kontrol af ansøgere til statsborgerskab og politiske holdninger. Flere parter udtrykker bekymringer over håndteringen af emner såsom kriminalitet og integration, samt behovet for klarhed og ansvarlighed i lovgivningen. Der er også et fokus på samarbejde mellem partierne i den politiske debat.'
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 'Der er en stigende politisk debat om indfødsret i Danmark, med fokus på at sikre, at statsborgerskab kun tildeles personer, der har vist loyalitet og tilknytning til landet. Der er bekymring over tildelingen af statsborgerskab til personer med kriminel baggrund og et ønske om strammere kriterier og individuel vurdering af ansøgere. Partier som Dansk Folkeparti presser på for at ændre procedurerne, hvilket skaber en kompleks diskussion om balance mellem integration, sikkerhed og danske værdier.'
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 'Der har været en række afstemninger om forskellige ændringsforslag, hvor mange blev forkastet med stor margin. Lovforslaget er blevet vedtaget til direkte behandling uden fornyet udvalgsbehandling, og der er flere ændringsforslag, der er blevet betragtet som vedtaget uden afstemning, når ingen har gjort indsigelse. Generelt set har ændringsforslag fra mindretal, især fra DF, ofte ikke fået støtte fra flertallet.'
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 'Der er en fælles opfattelse af, at Socialdemokratiet lider under utilstrækkelig intern kommunikation.'
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = 'Der er enighed blandt partierne om at stramme reglerne for dansk statsborgerskab, herunder krav om, at ansøgere skal dele grundlæggende demokratiske værdier og udelukkelse af personer med kriminel baggrund. Liberal Alliance og andre partier understreger behovet for at håndhæve disse stramninger og prioritere kvalificerede ansøgere.'
$ws.Range("A13").Value = 13
$ws.Range("B13").Value = 'Der er en igangværende debat om statsborgerskabslovgivning, hvor forskellige partier stiller spørgsmålstegn ved hinandens holdninger og forslag, herunder håndtering af statsløse personer og krav til ansøgere om statsborgerskab. Der er bekymring for, at massetildeling af statsborgerskaber uden tilstrækkelig kontrol kan have negative konsekvenser for samfundet, samtidig med at der er behov for at balancere demokratiske værdier med sikkerhed og integration. Diskussionen omfatter også hvordan kriminalitet påvirker tildeling af statsborgerskab og hvordan forskellige partier forholder sig til dette emne.'
$ws.Range("A14").Value = 14
$ws.Range("B14").Value = 'To personer med kriminel baggrund er blevet forpligtet til at optages på statsborgerskabslovforslaget på grund af sagsbehandlingsfejl, hvilket rejser spørgsmål om lovens anvendelse.'
$ws.Range("A15").Value = 15
$ws.Range("B15").Value = 'Der er en udbredt bekymring over indvandring og tildeling af dansk statsborgerskab, især til personer med kriminel baggrund, hvilket fører til uenighed mellem politiske partier om stramning af reglerne. Mange parter ønsker at sikre, at statsborgerskab kun gives til dem, der bidrager positivt til samfundet, og der efterlyses bedre kontrol med ansøgeres baggrund og rettigheder. Diskussionen om indvandringens omfang og dens indvirkning på det danske samfund er central, med krav om klare kriterier for tildeling af statsborgerskab.'
$ws.Range("A16").Value = 16
$ws.Range("B16").Value = 'Forhandlingen er afsluttet, da ingen ønsker at udtale sig, og der er indkaldt til afstemning. Der er stillet ændringsforslag, men ingen har ønsket at kommentere dem.'

$ws.Range("A17").Value = 17
$ws.Range("B17").Value = 'Folketinget står over for en debat om indfødsret, hvor der stilles spørgsmål ved, om internationale konventioner skal have indflydelse på beslutninger om statsborgerskab, især i forhold til personer, der kan udgøre en sikkerhedsrisiko. Der er også bekymringer om, hvordan sager behandles og beslutninger træffes i Indfødsretsudvalget, samt om det er acceptabelt at inkludere personer med problematiske baggrunde i lovforslag. Der er enighed om, at grundloven har forrang over konventionerne, og Folketinget har ret til at fastlægge betingelserne for indfødsret.'
$ws.Range("A18").Value = 18
$ws.Range("B18").Value = 'Der er en generel enighed blandt partierne om behovet for strammere kriterier for tildeling af dansk statsborgerskab, især i forhold til personer med kriminel baggrund. Der stilles spørgsmål ved, hvordan indfødsretsprocessen skal håndteres, herunder behovet for individuel behandling og objektive kriterier, samt bekymringer om politisk indflydelse på beslutningerne. Derudover er der en åbenhed for dialog og forhandlinger mellem partierne, selvom der er uenigheder om specifikke tilgange og holdninger til statsborgerskab.'
$ws.Range("A19").Value = 19
$ws.Range("B19").Value = 'Der rejses kritik af, at partier som Alternativet og Radikale Venstre ikke deltager i behandlingen af et vigtigt forslag, hvilket stiller spørgsmål ved deres seriøsitet. Liberal Alliance repræsenterer en ændring i tilgang ved at fokusere på ansøgeres værdier, hvilket skaber debat om partiets tidligere holdninger. Forhandlingen om lovforslaget er afsluttet, og det er vedtaget at henvise det til Indfødsretsudvalget.'

Write-Output "applied edits"